$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'327.99"
$ws.Range("E2").Value = "'-0.14%"
$ws.Range("G2").Value = "'17"
$ws.Range("D3").Value = "'44.44"
$ws.Range("E3").Value = "'0.69%"
$ws.Range("G3").Value = "'17"
$ws.Range("D4").Value = "'5.562"
$ws.Range("E4").Value = "'1.55%"
$ws.Range("G4").Value = "'17"
$ws.Range("D5").Value = "'0.08054"
$ws.Range("E5").Value = "'-1.29%"
$ws.Range("G5").Value = "'17"
$ws.Range("D6").Value = "'1.921"
$ws.Range("E6").Value = "'1.06%"
$ws.Range("G6").Value = "'17"
$ws.Range("D7").Value = "'2.574"
$ws.Range("E7").Value = "'-9.15%"
$ws.Range("G7").Value = "'17"
$ws.Range("D8").Value = "'0.9504"
$ws.Range("E8").Value = "'0.78%"
$ws.Range("G8").Value = "'17"
$ws.Range("D9").Value = "'0.1205"
$ws.Range("E9").Value = "'1.67%"
$ws.Range("G9").Value = "'17"
$ws.Range("D10").Value = "'0.1843"
$ws.Range("E10").Value = "'-3.07%"
$ws.Range("G10").Value = "'17"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.09687"
$ws.Range("E11").Value = "'-1.80%"
$ws.Range("G11").Value = "'17"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.04364"
$ws.Range("E12").Value = "'4.34%"
$ws.Range("G12").Value = "'17"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.1065"
$ws.Range("E13").Value = "'-0.25%"
$ws.Range("G13").Value = "'17"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001275"
$ws.Range("E14").Value = "'-2.40%"
$ws.Range("G14").Value = "'17"
$ws.Range("B15").Value = "CoinExToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D15").Value = "'0.04209"
$ws.Range("E15").Value = "'-3.98%"
$ws.Range("G15").Value = "'17"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005937"
$ws.Range("E16").Value = "'-2.50%"
$ws.Range("G16").Value = "'17"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.396"
$ws.Range("E17").Value = "'-3.94%"
$ws.Range("G17").Value = "'17"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "'4.281"
$ws.Range("E18").Value = "'-1.15%"
$ws.Range("G18").Value = "'17"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "'0.3464"
$ws.Range("E19").Value = "'-1.52%"
$ws.Range("G19").Value = "'17"
$ws.Range("B20").Value = "MCDex"
$ws.Range("C20").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D20").Value = "'10.01"
$ws.Range("E20").Value = "'13.74%"
$ws.Range("G20").Value = "'17"
$ws.Range("E21").Value = "'5.16%"
$ws.Range("G21").Value = "'17"
$ws.Range("D22").Value = "'0.2502"
$ws.Range("E22").Value = "'0.02%"
$ws.Range("G22").Value = "'17"
$ws.Range("D23").Value = "'0.001243"
$ws.Range("E23").Value = "'0.09%"
$ws.Range("G23").Value = "'17"
$ws.Range("D24").Value = "'0.004334"
$ws.Range("E24").Value = "'0.51%"
$ws.Range("G24").Value = "'17"
$ws.Range("E25").Value = "'-3.66%"
$ws.Range("G25").Value = "'17"
$ws.Range("E26").Value = "'-0.88%"
$ws.Range("G26").Value = "'17"
$ws.Range("G27").Value = "'17"
$ws.Range("G28").Value = "'17"
$ws.Range("G29").Value = "'17"
$ws.Range("G30").Value = "'17"
$ws.Range("G31").Value = "'17"
$ws.Range("G32").Value = "'17"
$ws.Range("G33").Value = "'17"
$ws.Range("G34").Value = "'17"
$ws.Range("G35").Value = "'17"
$ws.Range("G36").Value = "'17"
$ws.Range("G37").Value = "'17"
$ws.Range("D38").Value = "'0.02681"
$ws.Range("E38").Value = "'0.39%"
$ws.Range("G38").Value = "'17"
$ws.Range("D39").Value = "'0.05518"
$ws.Range("E39").Value = "'-2.47%"
$ws.Range("G39").Value = "'17"
$ws.Range("D40").Value = "'0.007548"
$ws.Range("E40").Value = "'-4.22%"
$ws.Range("G40").Value = "'17"
$ws.Range("D41").Value = "'0.1405"
$ws.Range("E41").Value = "'-0.30%"
$ws.Range("G41").Value = "'17"
$ws.Range("D42").Value = "'0.008303"
$ws.Range("E42").Value = "'-14.82%"
$ws.Range("G42").Value = "'17"
$ws.Range("D43").Value = "'0.002017"
$ws.Range("E43").Value = "'-4.38%"
$ws.Range("G43").Value = "'17"
$ws.Range("D44").Value = "'0.008892"
$ws.Range("E44").Value = "'-7.92%"
$ws.Range("G44").Value = "'17"
$ws.Range("D45").Value = "'0.00007107"
$ws.Range("E45").Value = "'0.59%"
$ws.Range("G45").Value = "'17"
$ws.Range("D46").Value = "'0.00000000749"
$ws.Range("E46").Value = "'-0.63%"
$ws.Range("G46").Value = "'17"
$ws.Range("D47").Value = "'0.002829"
$ws.Range("E47").Value = "'-18.14%"
$ws.Range("G47").Value = "'17"
$ws.Range("D48").Value = "'0.002267"
$ws.Range("E48").Value = "'-0.75%"
$ws.Range("G48").Value = "'17"
$ws.Range("D49").Value = "'0.00002098"
$ws.Range("E49").Value = "'-0.63%"
$ws.Range("G49").Value = "'17"
$ws.Range("D50").Value = "'0.0001998"
$ws.Range("E50").Value = "'-0.63%"
$ws.Range("G50").Value = "'17"
$ws.Range("G51").Value = "'17"

Write-Output "Applied 143 changes"